$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)

# Add a paragraph border (top/left/bottom/right, 5pt space from text) and
# change the left indent from 120 twips (6pt) to 225 twips (11.25pt).
$b = $p1.Format.Borders
$b.DistanceFromTop = 5
$b.DistanceFromBottom = 5
$b.DistanceFromLeft = 5
$b.DistanceFromRight = 5
$p1.Format.LeftIndent = 11.25

# Replace the paragraph's text (which spans two runs: the id placeholder and
# a trailing space run) with the updated placeholder text and no trailing
# space, collapsing it into a single run.
$r = $p1.Range
$body = $d.Range($r.Start, $r.End - 1)
$body.Text = "**ID__AFFARS_5322_7003__ID**"
